$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (columns B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values - C2 and E2 are cleared (removed), B2 and D2 updated
$ws.Range("B2").Value = 23.254962237594334
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 35.204401339886694
$ws.Range("E2").ClearContents()

# Row 3 values updated
$ws.Range("B3").Value = 20.682618561610106
$ws.Range("C3").Value = -6.5016201590062561
$ws.Range("D3").Value = 31.903132892840148
$ws.Range("E3").Value = -0.79204125872306064

# Update the selection to match new sqref B1:E3
$ws.Range("B1:E3").Select() | Out-Null
